$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values as text (e.g. "0.5200", "1.270", "26.080.78").
# Force the Price column to a Text number format before writing the new
# values so Excel does not auto-convert numeric-looking strings into true
# numbers (which would silently drop significant trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.080.78'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.651.17'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '218.24'
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').Value = '0.5200'
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').Value = '0.2642'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').Value = '0.06330'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('D10').Value = '20.36'
$ws.Range('E10').Value = '  -0.52%  '
$ws.Range('D11').Value = '0.07682'
$ws.Range('E11').Value = '  -1.41%  '
$ws.Range('D12').Value = '4.593'
$ws.Range('E12').Value = '  +2.67%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.652.74'
$ws.Range('E13').Value = '  +3.45%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.880.40'
$ws.Range('E14').Value = '  +0.20%  '
$ws.Range('D15').Value = '0.5583'
$ws.Range('E15').Value = '  +1.06%  '
$ws.Range('D16').Value = '0.0₅8137'
$ws.Range('E16').Value = '  +1.94%  '
$ws.Range('D17').Value = '65.33'
$ws.Range('E17').Value = '  +0.97%  '
$ws.Range('D18').Value = '26.084.93'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('E19').Value = '  -0.29%  '
$ws.Range('D20').Value = '4.624'
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('D21').Value = '10.47'
$ws.Range('E21').Value = '  +4.20%  '
$ws.Range('D22').Value = '191.29'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').Value = '5.914'
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('D25').Value = '144.23'
$ws.Range('E25').Value = '  -1.70%  '
$ws.Range('D26').Value = '0.1186'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('D27').Value = '7.215'
$ws.Range('E27').Value = '  +0.72%  '
$ws.Range('D28').Value = '15.89'
$ws.Range('E28').Value = '  +0.18%  '
$ws.Range('D29').Value = '1.501'
$ws.Range('E29').Value = '  +1.60%  '
$ws.Range('D30').Value = '0.05477'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('D31').Value = '1.270'
$ws.Range('E31').Value = '  +0.51%  '
$ws.Range('D32').Value = '3.441'
$ws.Range('D33').Value = '3.351'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').Value = '1.557'
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').Value = '2.425'
$ws.Range('E35').Value = '  +0.48%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').Value = '2.786'
$ws.Range('E36').Value = '  -0.32%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').Value = '0.9454'
$ws.Range('E37').Value = '  -0.12%  '
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('D39').Value = '0.01581'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').Value = '5.850'
$ws.Range('E40').Value = '  -1.68%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('D42').Value = '1.026.17'
$ws.Range('E42').Value = '  -2.95%  '
$ws.Range('D43').Value = '0.8256'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('D44').Value = '101.17'
$ws.Range('E44').Value = '  -0.95%  '
$ws.Range('D45').Value = '1.794.77'
$ws.Range('E45').Value = '  +0.39%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈109'
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '57.49'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').Value = '0.9992'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.998'
$ws.Range('E49').Value = '  +1.21%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '0.4336'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('D51').Value = '0.05165'
$ws.Range('E51').Value = '  -3.39%  '
